# edit.ps1 - applies the "Report" docx revision:
#   1. Paragraph line spacing: 480 (auto) -> 600 (auto), i.e. 24pt -> 30pt multiple.
#   2. Font size for every run (and paragraph mark) in both paragraphs: 24 half-pts (12pt) -> 28 half-pts (14pt).
#   3. Appends additional sentence text (as four separate runs, matching the
#      original author's edit history) to the end of the second paragraph,
#      just before the trailing "_GoBack" bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1 & 2: spacing + font size for both paragraphs in the document.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 5      # wdLineSpaceMultiple
    $p.Format.LineSpacing = 30         # 30pt ~ "line 600" (600/20)

    $p.Range.Font.Size = 14            # 28 half-points
    $p.Range.Font.SizeBi = 14          # keep complex-script size (szCs) in sync
}

# ---------------------------------------------------------------------------
# 3: append the new sentence to paragraph 2, split across 4 runs exactly as
#    in the authored edit, inserted before the "_GoBack" bookmark so that the
#    bookmark remains collapsed at the very end of the paragraph afterward.
# ---------------------------------------------------------------------------

# Remove the existing _GoBack bookmark (Word re-creates it automatically as
# the user edits; we'll recreate a collapsed one in the correct spot below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the end of the existing text in paragraph 2 ("...left as it is.")
$insertPoint = $d.Content
$insertPoint.Find.Execute("left as it is.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint.Collapse(0)   # wdCollapseEnd

$newRuns = @(
    " Finally, the output prefix was a single string so we had to divide the string into an array of strings so that we could deal with every element in i",
    "ts",
    " own",
    "."
)

foreach ($t in $newRuns) {
    $insertPoint.InsertAfter($t)
    $insertPoint.Collapse(0)   # wdCollapseEnd - move past what we just inserted
    $insertPoint.Font.Size = 14
    $insertPoint.Font.SizeBi = 14
}

# Recreate a collapsed "_GoBack" bookmark at the true end of the paragraph's
# text (mirrors how Word leaves the bookmark after the last edit position).
$insertPoint.InsertAfter("X")
$placeholder = $d.Range($insertPoint.Start, $insertPoint.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$d.Bookmarks("_GoBack").Range.Delete()
